# Class 5 - Data Types in Java
# Fix the "short" data type value-range line on the "Short Data type" slide:
#   "-2^16 to +2^16-1"  ->  "-2^15 to +2^15-1"
# split across three runs ("-2^15 ", "to ", "+2^15-1") matching the way the
# author actually retyped the first and last numbers while leaving the
# middle "to " text untouched.

$p = $ppt.ActivePresentation

$oldText = "-2^16 to +2^16-1"
$newFirst = "-2^15 "
$newLast = "+2^15-1"

$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf($oldText)
        if ($idx -ge 0) {
            $start = $idx + 1

            # First segment: "-2^16 " (6 chars) -> "-2^15 "
            $seg1 = $tr.Characters($start, 6)
            $seg1.Text = $newFirst

            # Last segment: "+2^16-1" (7 chars), starts 9 chars after $start
            # (unaffected by the seg1 edit above since the replacement text
            # is the same length).
            $seg3 = $tr.Characters($start + 9, 7)
            $seg3.Text = $newLast

            $found = $true
        }
    }
}

if (-not $found) {
    throw "Could not locate the '-2^16 to +2^16-1' text to update."
}
